# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 4 of the sheet
# (pushing the previously existing rows 4-29 down to rows 5-30).
# The new row carries a new "Fecha" (date serial 44473) together with
# updated Volumen / Precio / Origen figures, while every other column
# repeats the constant values used throughout this "Haba" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (rows 4:29) down one row, just like using
# Excel's "Insert Sheet Rows" above the current row 4.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's record.
$ws.Cells.Item(4, 1).Value  = 10
$ws.Cells.Item(4, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value  = "La Araucanía"
$ws.Cells.Item(4, 4).Value  = 44473
$ws.Cells.Item(4, 5).Value  = 9
$ws.Cells.Item(4, 6).Value  = 100112026
$ws.Cells.Item(4, 7).Value  = "Haba"
$ws.Cells.Item(4, 8).Value  = "Sin especificar"
$ws.Cells.Item(4, 9).Value  = "Primera"
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 13).Value = 10000
$ws.Cells.Item(4, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 400
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
